{"js": "// Update the division problems in the worksheet table.\n// The table has 20 rows x 5 columns; only rows 0, 4, 8, 12, 16 (0-indexed)\n// contain text. Each entry below is [rowIndex, colIndex, oldText, newText].\nconst replacements = [\n  [0, 0, \"41\u00f76=\", \"25\u00f72=\"],\n  [0, 1, \"12\u00f78=\", \"68\u00f79=\"],\n  [0, 2, \"39\u00f74=\", \"27\u00f72=\"],\n  [0, 3, \"12\u00f72=\", \"90\u00f72=\"],\n  [0, 4, \"70\u00f77=\", \"89\u00f76=\"],\n  [4, 0, \"54\u00f79=\", \"99\u00f73=\"],\n  [4, 1, \"20\u00f72=\", \"44\u00f79=\"],\n  [4, 2, \"80\u00f79=\", \"33\u00f76=\"],\n  [4, 3, \"79\u00f73=\", \"30\u00f75=\"],\n  [4, 4, \"99\u00f79=\", \"25\u00f77=\"],\n  [8, 0, \"73\u00f72=\", \"42\u00f75=\"],\n  [8, 1, \"94\u00f77=\", \"90\u00f72=\"],\n  [8, 2, \"72\u00f78=\", \"88\u00f75=\"],\n  [8, 3, \"52\u00f75=\", \"23\u00f76=\"],\n  [8, 4, \"83\u00f78=\", \"92\u00f73=\"],\n  [12, 0, \"79\u00f73=\", \"42\u00f78=\"],\n  [12, 1, \"80\u00f76=\", \"70\u00f79=\"],\n  [12, 2, \"53\u00f75=\", \"90\u00f74=\"],\n  [12, 3, \"38\u00f74=\", \"20\u00f72=\"],\n  [12, 4, \"21\u00f74=\", \"87\u00f79=\"],\n  [16, 0, \"12\u00f76=\", \"89\u00f76=\"],\n  [16, 1, \"52\u00f78=\", \"82\u00f74=\"],\n  [16, 2, \"85\u00f77=\", \"57\u00f74=\"],\n  [16, 3, \"79\u00f72=\", \"86\u00f72=\"],\n  [16, 4, \"31\u00f75=\", \"50\u00f75=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [rowIndex, colIndex, oldText, newText] of replacements) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\n      `Expected text \"${oldText}\" not found in cell (${rowIndex}, ${colIndex})`\n    );\n  }\n  // Replace just the matched range so run/paragraph formatting is preserved.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the worksheet table.\n# The table has 20 rows x 5 columns (1-indexed via COM); only rows\n# 1, 5, 9, 13, 17 contain text. Each entry below identifies a cell by\n# (Row, Col) together with the expected current text and its replacement,\n# so cells are addressed positionally and duplicate text values (e.g.\n# \"79\u00f73=\" appears twice, mapping to two different results) are handled\n# unambiguously.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"41\u00f76=\"; NewText = \"25\u00f72=\" },\n    @{ Row = 1; Col = 2; OldText = \"12\u00f78=\"; NewText = \"68\u00f79=\" },\n    @{ Row = 1; Col = 3; OldText = \"39\u00f74=\"; NewText = \"27\u00f72=\" },\n    @{ Row = 1; Col = 4; OldText = \"12\u00f72=\"; NewText = \"90\u00f72=\" },\n    @{ Row = 1; Col = 5; OldText = \"70\u00f77=\"; NewText = \"89\u00f76=\" },\n    @{ Row = 5; Col = 1; OldText = \"54\u00f79=\"; NewText = \"99\u00f73=\" },\n    @{ Row = 5; Col = 2; OldText = \"20\u00f72=\"; NewText = \"44\u00f79=\" },\n    @{ Row = 5; Col = 3; OldText = \"80\u00f79=\"; NewText = \"33\u00f76=\" },\n    @{ Row = 5; Col = 4; OldText = \"79\u00f73=\"; NewText = \"30\u00f75=\" },\n    @{ Row = 5; Col = 5; OldText = \"99\u00f79=\"; NewText = \"25\u00f77=\" },\n    @{ Row = 9; Col = 1; OldText = \"73\u00f72=\"; NewText = \"42\u00f75=\" },\n    @{ Row = 9; Col = 2; OldText = \"94\u00f77=\"; NewText = \"90\u00f72=\" },\n    @{ Row = 9; Col = 3; OldText = \"72\u00f78=\"; NewText = \"88\u00f75=\" },\n    @{ Row = 9; Col = 4; OldText = \"52\u00f75=\"; NewText = \"23\u00f76=\" },\n    @{ Row = 9; Col = 5; OldText = \"83\u00f78=\"; NewText = \"92\u00f73=\" },\n    @{ Row = 13; Col = 1; OldText = \"79\u00f73=\"; NewText = \"42\u00f78=\" },\n    @{ Row = 13; Col = 2; OldText = \"80\u00f76=\"; NewText = \"70\u00f79=\" },\n    @{ Row = 13; Col = 3; OldText = \"53\u00f75=\"; NewText = \"90\u00f74=\" },\n    @{ Row = 13; Col = 4; OldText = \"38\u00f74=\"; NewText = \"20\u00f72=\" },\n    @{ Row = 13; Col = 5; OldText = \"21\u00f74=\"; NewText = \"87\u00f79=\" },\n    @{ Row = 17; Col = 1; OldText = \"12\u00f76=\"; NewText = \"89\u00f76=\" },\n    @{ Row = 17; Col = 2; OldText = \"52\u00f78=\"; NewText = \"82\u00f74=\" },\n    @{ Row = 17; Col = 3; OldText = \"85\u00f77=\"; NewText = \"57\u00f74=\" },\n    @{ Row = 17; Col = 4; OldText = \"79\u00f72=\"; NewText = \"86\u00f72=\" },\n    @{ Row = 17; Col = 5; OldText = \"31\u00f75=\"; NewText = \"50\u00f75=\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $rng = $cell.Range\n    # The cell range's text includes trailing cell-mark control\n    # characters (cr/bell); strip those off before comparing so the\n    # sanity check is exact.\n    $current = $rng.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $item.OldText) {\n        throw \"Cell ($($item.Row),$($item.Col)) expected '$($item.OldText)' but found '$current'\"\n    }\n    # Assigning straight to the cell's Range.Text replaces the visible\n    # run content while Word keeps the cell-end mark and reuses the\n    # existing run/paragraph formatting (font, size, alignment) instead\n    # of resetting it.\n    $rng.Text = $item.NewText\n}\n"}
